$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (col A and col E changed) ---
$ws.Columns.Item(1).ColumnWidth = 72.33203125
$ws.Columns.Item(5).ColumnWidth = 30.6640625

# --- Clear the stray wrap-text style that used to sit on E4 (no longer used) ---
$ws.Range("E4").WrapText = $false

# --- Row 2 ---
$ws.Range("A2").Value = "LD_D200_Contract_Financial_JobSeq"
$ws.Range("B2").Value = "Sequence"
$ws.Range("E2").Value = "D200 CNTRCTFI-BMSIW ETL 01 D"

# --- Row 3 ---
$ws.Range("A3").Value = "LD_D200_REVENUE_COST_CATEGORY_REFERENCE_PJob"
$ws.Range("B3").Value = "Parallel"
$ws.Range("C3").Value = "J200104"
$ws.Range("D3").Value = "C200104"
$ws.Range("E3").Value = ""
$ws.Range("E3").WrapText = $true
$ws.Range("F3").Value = "LD_D200_Contract_Financial_JobSeq"

# --- Row 4 ---
$ws.Range("A4").Value = "LD_D200_LAST_YEAR_REVENUE_COST_HISTORY_PJob"
$ws.Range("B4").Value = "Parallel"
$ws.Range("C4").Value = "J200105"
$ws.Range("D4").Value = "C200105"
$ws.Range("F4").Value = "LD_D200_Contract_Financial_JobSeq"

# --- Row 5 ---
$ws.Range("A5").Value = "LD_D200_REVENUE_COST_CURRENT_HISTORY_PJob"
$ws.Range("B5").Value = "Parallel"
$ws.Range("C5").Value = "J200106"
$ws.Range("D5").Value = "C200106"
$ws.Range("F5").Value = "LD_D200_Contract_Financial_JobSeq"

# --- Row 6 ---
$ws.Range("A6").Value = "LD_D200_CURRENT_MONTH_REVENUE_DETAIL_COST_PJob"
$ws.Range("B6").Value = "Parallel"
$ws.Range("C6").Value = "J200108"
$ws.Range("D6").Value = "C200108"
$ws.Range("F6").Value = "LD_D200_Contract_Financial_JobSeq"

# --- Row 7 ---
$ws.Range("A7").Value = "LD_D200_CURRENT_YEAR_REVENUE_DETAIL_COST_PJob"
$ws.Range("B7").Value = "Parallel"
$ws.Range("C7").Value = "J200109"
$ws.Range("D7").Value = "C200109"
$ws.Range("F7").Value = "LD_D200_Contract_Financial_JobSeq"

# --- Row 8 ---
$ws.Range("A8").Value = "LD_D200_WORK_NUMBER_REVENUE_REFERENCE_PJob"
$ws.Range("B8").Value = "Parallel"
$ws.Range("C8").Value = "J200110"
$ws.Range("D8").Value = "C200110"
$ws.Range("F8").Value = "LD_D200_Contract_Financial_JobSeq"

# --- Hyperlinks: drop the old two and recreate the two required by the new data ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "http://iwrefresh.w3ibm.mybluemix.net/Domains/ODS ADL/Datagroups/D800 CUSTOMER-BMSIW ETL 01 D", "", "", "http://iwrefresh.w3ibm.mybluemix.net/Domains/ODS ADL/Datagroups/D800 CUSTOMER-BMSIW ETL 01 D")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:LD_D200_REVENUE_COST_CATEGORY_REFERENCE_PJob@J200104@C200104", "", "", "LD_D200_REVENUE_COST_CATEGORY_REFERENCE_PJob@J200104@C200104")

# --- Selection matches the saved cursor position in the new workbook ---
$ws.Range("A12").Select()
